# RequirementsTracing.xlsx - add the "advanced booking" requirements-tracing
# row content (requirement 3.1.9 / advancedBooking) and refresh the view
# metadata that Excel rewrote when the sheet was resaved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core content edit: row 11 (3.1.9) gained the advancedBooking() tracing ---
$ws.Range("E11").Value = "public static boolean advancedBooking() throws Exception"
$ws.Range("F11").Value = "Main"
$ws.Range("G11").Value = "Line 243 - 285"
$ws.Range("H11").Value = "When an advanced booking is made, the operator must enter the customers' details into the system. The booking is then made using the customers details. If the rating of the showing is greater than the customers' age, the booking is rejected."

# Those two cells wrap their (now longer) text, same as the other "Method
# Description" / long-text cells in the table.
$ws.Range("E11").WrapText = $true
$ws.Range("H11").WrapText = $true

# Row grew taller to fit the new wrapped text.
$ws.Rows(11).RowHeight = 105

# --- View state: scrolled down a bit and selection moved to H12 ---
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("H12").Select()

# --- Column widths settled to their recalculated best-fit values ---
$ws.Columns("A").ColumnWidth = 2.75
$ws.Columns("B").ColumnWidth = 15.375
$ws.Columns("C").ColumnWidth = 35.125
$ws.Columns("D").ColumnWidth = 36.125
$ws.Columns("E").ColumnWidth = 36.625
$ws.Columns("F").ColumnWidth = 7.125
$ws.Columns("G").ColumnWidth = 13
$ws.Columns("H").ColumnWidth = 37.125
$ws.Columns("I").ColumnWidth = 5

# --- Row heights settled to their recalculated values ---
$ws.Rows(1).RowHeight = 13.5
$ws.Rows(3).RowHeight = 93.95
$ws.Rows(5).RowHeight = 204.95
$ws.Rows(6).RowHeight = 62.1
$ws.Rows(7).RowHeight = 167.1
$ws.Rows(8).RowHeight = 30.95
$ws.Rows(10).RowHeight = 84.95
$ws.Rows(12).RowHeight = 35.1
$ws.Rows(13).RowHeight = 32.1
$ws.Rows(14).RowHeight = 33.95
$ws.Rows(16).RowHeight = 18.75
For ($r = 17; $r -le 40; $r++) {
    $ws.Rows($r).RowHeight = 15
}
